# Auto-generated Excel COM-interop script applying the Bahamut_Profits leve-profit refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1806.5
$ws.Range("I135").Value = 1806.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 16258.5
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -13723.5
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 19800
$ws.Range("J29").Value = 19800
$ws.Range("L29").Value = 19800
$ws.Range("N29").Value = -20416

$ws.Range("H32").Value = 5023.77
$ws.Range("I32").Value = 4872.495
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 4872.495
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -4585.495
$ws.Range("N32").Value = -20574

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H74").Value = 1524
$ws.Range("I74").Value = 1548.9412
$ws.Range("J74").Value = 1100
$ws.Range("K74").Value = 1548.9412
$ws.Range("L74").Value = 1100
$ws.Range("M74").Value = -674.9412
$ws.Range("N74").Value = -2848

$ws.Range("H77").Value = 1524
$ws.Range("I77").Value = 1548.9412
$ws.Range("J77").Value = 1100
$ws.Range("K77").Value = 7744.706
$ws.Range("L77").Value = 5500
$ws.Range("M77").Value = -3376.706
$ws.Range("N77").Value = -14236

$ws.Range("H80").Value = 27444
$ws.Range("J80").Value = 27444
$ws.Range("L80").Value = 27444
$ws.Range("N80").Value = -29440

$ws.Range("H83").Value = 27444
$ws.Range("J83").Value = 27444
$ws.Range("L83").Value = 82332
$ws.Range("N83").Value = -92316

$ws.Range("H122").Value = 820
$ws.Range("J122").Value = 800
$ws.Range("L122").Value = 2400
$ws.Range("N122").Value = -7300

$ws.Range("H132").Value = 2023.12
$ws.Range("I132").Value = 1628.5883
$ws.Range("J132").Value = 2861.5
$ws.Range("K132").Value = 4885.7649
$ws.Range("L132").Value = 8584.5
$ws.Range("M132").Value = -2355.7649
$ws.Range("N132").Value = -13644.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 888.1667
$ws.Range("I94").Value = 907.2174
$ws.Range("K94").Value = 907.2174
$ws.Range("M94").Value = -456.2174

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2705.3242
$ws.Range("I31").Value = 3226.7896
$ws.Range("J31").Value = 2154.889
$ws.Range("K31").Value = 3226.7896
$ws.Range("L31").Value = 2154.889
$ws.Range("M31").Value = -2931.7896
$ws.Range("N31").Value = -2744.889

$ws.Range("H34").Value = 2705.3242
$ws.Range("I34").Value = 3226.7896
$ws.Range("J34").Value = 2154.889
$ws.Range("K34").Value = 3226.7896
$ws.Range("L34").Value = 2154.889
$ws.Range("M34").Value = -3024.7896
$ws.Range("N34").Value = -2558.889

$ws.Range("H58").Value = 3681.0889
$ws.Range("I58").Value = 1606.9565
$ws.Range("K58").Value = 1606.9565
$ws.Range("M58").Value = -1403.9565

$ws.Range("H127").Value = 54163.332
$ws.Range("J127").Value = 54163.332
$ws.Range("L127").Value = 54163.332
$ws.Range("N127").Value = -64083.332

$ws.Range("H134").Value = 16130438
$ws.Range("I134").Value = 1397.5
$ws.Range("J134").Value = 71430000
$ws.Range("K134").Value = 4192.5
$ws.Range("L134").Value = 214290000
$ws.Range("M134").Value = -1657.5
$ws.Range("N134").Value = -214295070

$ws.Range("H136").Value = 3681.0889
$ws.Range("I136").Value = 1606.9565
$ws.Range("K136").Value = 4820.8695
$ws.Range("M136").Value = -2270.8695

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 5000
$ws.Range("J31").Value = 5000
$ws.Range("L31").Value = 15000
$ws.Range("N31").Value = -15576

$ws.Range("H35").Value = 2500
$ws.Range("J35").Value = 2500
$ws.Range("L35").Value = 7500
$ws.Range("N35").Value = -8076

$ws.Range("H129").Value = 26909.975
$ws.Range("I129").Value = 1034.9166
$ws.Range("J129").Value = 37999.285
$ws.Range("K129").Value = 3104.7498
$ws.Range("L129").Value = 113997.855
$ws.Range("M129").Value = 1895.2502
$ws.Range("N129").Value = -123997.855

$ws.Range("H132").Value = 932.6087
$ws.Range("I132").Value = 666.3158
$ws.Range("J132").Value = 2197.5
$ws.Range("K132").Value = 5996.8422
$ws.Range("L132").Value = 19777.5
$ws.Range("M132").Value = -3466.8422
$ws.Range("N132").Value = -24837.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 13645.143
$ws.Range("I29").Value = 1500
$ws.Range("J29").Value = 18503.2
$ws.Range("K29").Value = 1500
$ws.Range("L29").Value = 18503.2
$ws.Range("M29").Value = -1210
$ws.Range("N29").Value = -19083.2

$ws.Range("H31").Value = 5696
$ws.Range("I31").Value = 870
$ws.Range("J31").Value = 25000
$ws.Range("K31").Value = 870
$ws.Range("L31").Value = 25000
$ws.Range("M31").Value = -578
$ws.Range("N31").Value = -25584

$ws.Range("H35").Value = 80017
$ws.Range("J35").Value = 80017
$ws.Range("L35").Value = 80017
$ws.Range("N35").Value = -80613

$ws.Range("H36").Value = 10182.357
$ws.Range("I36").Value = 17052.834
$ws.Range("J36").Value = 5029.5
$ws.Range("K36").Value = 17052.834
$ws.Range("L36").Value = 5029.5
$ws.Range("M36").Value = -16567.834
$ws.Range("N36").Value = -5999.5

$ws.Range("H37").Value = 5696
$ws.Range("I37").Value = 870
$ws.Range("J37").Value = 25000
$ws.Range("K37").Value = 870
$ws.Range("L37").Value = 25000
$ws.Range("M37").Value = -593
$ws.Range("N37").Value = -25554

$ws.Range("H70").Value = 5171.2856
$ws.Range("I70").Value = 5127.091
$ws.Range("K70").Value = 5127.091
$ws.Range("M70").Value = -4857.091

$ws.Range("H73").Value = 5171.2856
$ws.Range("I73").Value = 5127.091
$ws.Range("K73").Value = 5127.091
$ws.Range("M73").Value = -4191.091

$ws.Range("H122").Value = 823592.4399999999
$ws.Range("I122").Value = 941047.9399999999
$ws.Range("J122").Value = 1404
$ws.Range("K122").Value = 2823143.82
$ws.Range("L122").Value = 4212
$ws.Range("M122").Value = -2820693.82
$ws.Range("N122").Value = -9112

$ws.Range("H128").Value = 27666.555
$ws.Range("J128").Value = 27666.555
$ws.Range("L128").Value = 27666.555
$ws.Range("N128").Value = -37626.555

$ws.Range("H132").Value = 2394.5
$ws.Range("I132").Value = 2023.7812
$ws.Range("J132").Value = 3383.0833
$ws.Range("K132").Value = 6071.3436
$ws.Range("L132").Value = 10149.2499
$ws.Range("M132").Value = -3541.3436
$ws.Range("N132").Value = -15209.2499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H30").Value = 7809.091
$ws.Range("J30").Value = 7809.091
$ws.Range("L30").Value = 7809.091
$ws.Range("N30").Value = -8025.091

$ws.Range("H31").Value = 1447.5
$ws.Range("I31").Value = 1142.5
$ws.Range("J31").Value = 1600
$ws.Range("K31").Value = 1142.5
$ws.Range("L31").Value = 1600
$ws.Range("M31").Value = -894.5
$ws.Range("N31").Value = -2096

$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H122").Value = 2090.4546
$ws.Range("I122").Value = 2030.3125
$ws.Range("J122").Value = 2250.8333
$ws.Range("K122").Value = 6090.9375
$ws.Range("L122").Value = 6752.499899999999
$ws.Range("M122").Value = -3640.9375
$ws.Range("N122").Value = -11652.4999

$ws.Range("H132").Value = 1732.1464
$ws.Range("I132").Value = 1500.5264
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 4501.5792
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -1971.5792
$ws.Range("N132").Value = -19058

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

$ws.Range("H132").Value = 1654.4445
$ws.Range("I132").Value = 1634.1818
$ws.Range("J132").Value = 1686.2858
$ws.Range("K132").Value = 4902.5454
$ws.Range("L132").Value = 5058.857400000001
$ws.Range("M132").Value = -2372.5454
$ws.Range("N132").Value = -10118.8574
